$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data cells are text-formatted (prices/percentages stored as
# strings like "10.90" or "1.004"). Force NumberFormat to Text ("@") before
# assignment so Excel does not auto-convert them to numbers and strip
# formatting (e.g. trailing zeros).

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.751.42'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +6.34%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.736.83'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +5.11%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.005'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.17%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '227.53'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +4.08%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5447'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +3.68%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.005'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.11%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2771'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +3.77%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06728'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +5.66%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.89'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +6.54%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07778'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.20%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.691'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +2.10%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.777.44'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +3.44%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.978.40'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +5.24%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5977'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +6.57%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₅8392'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.00%  '

# Row 17
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +5.55%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '27.789.71'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +6.47%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '224.47'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +17.45%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.835'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +3.09%  '

# Row 21
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.06%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.90'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +5.29%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.233'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +4.38%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.006'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.13%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.43'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.27%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1247'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +3.90%  '

# Row 27
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +11.99%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.458'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +2.86%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '17.19'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +7.78%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05665'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.50%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.312'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +3.20%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.690'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +5.49%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.518'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +4.12%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.684'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +6.65%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9761'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +3.26%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.859'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +2.19%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.453'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +1.88%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5960'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +3.19%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01665'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +4.53%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.991'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.38%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8508'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.42%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.045.71'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +2.26%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.005'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.11%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.88'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.35%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.885.36'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +5.25%  '

# Row 46
$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0₈115'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +11.81%  '

# Row 47
$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '59.45'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.49%  '

# Row 48
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.259'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +2.97%  '

# Row 49
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4443'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +2.28%  '

# Row 50
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'Frax'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.004'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.08%  '

# Row 51
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05320'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.45%  '
